# Update NATMI ligand-receptor pair stats (Vegfb-Nrp1) after recomputation
# with revised per-cluster cell counts (Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T for data rows 2-37
$data = @{
    2 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 100.1880891077553; R = 901.6928019697981; S = 0.01575565309554173; T = 0.01575565309554173 }
    3 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 35.90489792938267; R = 323.144081364444; S = 0.005646430840673632; T = 0.005646430840673631 }
    4 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 51.86009849815934; R = 466.740886483434; S = 0.008155557499043804; T = 0.008155557499043804 }
    5 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 65.69929852406601; R = 591.293686716594; S = 0.01033192034486556; T = 0.01033192034486555 }
    6 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 20.16654798937734; R = 181.498931904396; S = 0.00317140627279045; T = 0.003171406272790449 }
    7 = @{ E = 3.0; G = 1.028415333333333; H = 3.085246; I = 0.04565156193945813; J = 0.04565156193945813; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 16.47323976817022; R = 148.259157913532; S = 0.002590593886542959; T = 0.002590593886542959 }
    8 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 479.8523232696004; R = 4318.670909426403; S = 0.07546193175113004; T = 0.07546193175113003 }
    9 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 171.9670356187927; R = 1547.703320569134; S = 0.02704366338561729; T = 0.02704366338561728 }
    10 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 248.3847029217943; R = 2235.462326296149; S = 0.03906116234301996; T = 0.03906116234301996 }
    11 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 314.667754567601; R = 2832.009791108409; S = 0.04948488413615642; T = 0.04948488413615641 }
    12 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 96.58797758506735; R = 869.2917982656062; S = 0.01518949689112777; T = 0.01518949689112776 }
    13 = @{ E = 3.0; G = 4.925610333333334; H = 14.776831; I = 0.2186488259495045; J = 0.2186488259495045; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 78.89882365190023; R = 710.089412867102; S = 0.01240768744245305; T = 0.01240768744245304 }
    14 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 426.056161251271; R = 3834.505451261439; S = 0.06700190746899465; T = 0.06700190746899462 }
    15 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 152.687840622038; R = 1374.190565598342; S = 0.02401180290164827; T = 0.02401180290164827 }
    16 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 220.538336293393; R = 1984.845026640537; S = 0.03468202210314089; T = 0.03468202210314088 }
    17 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 279.390406338213; R = 2514.513657043917; S = 0.04393714222608703; T = 0.04393714222608702 }
    18 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 85.759515912142; R = 771.8356432092781; S = 0.01348660498854357; T = 0.01348660498854356 }
    19 = @{ E = 3.0; G = 4.373401; H = 13.120203; I = 0.1941361434105301; J = 0.1941361434105301; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 70.05349000568067; R = 630.481410051126; S = 0.01101666372211571; T = 0.01101666372211571 }
    20 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 557.7382436757773; R = 5019.644193081996; S = 0.08771032927897258; T = 0.08771032927897257 }
    21 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 199.8793957328987; R = 1798.914561596088; S = 0.03143318180993512; T = 0.03143318180993512 }
    22 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 288.7005881717853; R = 2598.305293546068; S = 0.04540126831664867; T = 0.04540126831664867 }
    23 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 365.742192468932; R = 3291.679732220388; S = 0.057516888067858; T = 0.05751688806785799 }
    24 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 112.2653915926213; R = 1010.388524333592; S = 0.01765493862914258; T = 0.01765493862914258 }
    25 = @{ E = 3.0; G = 5.725097333333333; H = 17.175292; I = 0.2541382134735057; J = 0.2541382134735057; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 91.70507090985156; R = 825.345638188664; S = 0.01442160737094877; T = 0.01442160737094877 }
    26 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 175.0022960347217; R = 1575.020664312495; S = 0.0275209906866358; T = 0.0275209906866358 }
    27 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 62.71643298612334; R = 564.4478968751101; S = 0.00986283270116447; T = 0.009862832701164472 }
    28 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 90.58598073473166; R = 815.273826612585; S = 0.01424561842117596; T = 0.01424561842117596 }
    29 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 114.759430906165; R = 1032.834878155485; S = 0.01804715309875457; T = 0.01804715309875457 }
    30 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 35.22566637077667; R = 317.0309973369901; S = 0.005539614382706956; T = 0.005539614382706956 }
    31 = @{ E = 3.0; G = 1.796371666666667; H = 5.389115; I = 0.07974129687595831; J = 0.07974129687595832; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 28.77442626398111; R = 258.96983637583; S = 0.004525087585520558; T = 0.004525087585520559 }
    32 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 97.419871; N = 292.259613; O = 0.3451284562056485; P = 0.3451284562056485; Q = 455.7885435140194; R = 4102.096891626174; S = 0.07167764392437373; T = 0.07167764392437372 }
    33 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 34.912838; N = 104.738514; O = 0.123685381195977; P = 0.123685381195977; Q = 163.3431805915747; R = 1470.088625324172; S = 0.02568746955693818; T = 0.02568746955693817 }
    34 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 50.427193; N = 151.281579; O = 0.1786479400170247; P = 0.1786479400170247; Q = 235.9286315516713; R = 2123.357683965042; S = 0.03710231133399541; T = 0.03710231133399541 }
    35 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 63.884013; N = 191.652039; O = 0.2263212890408322; P = 0.2263212890408321; Q = 298.887700632458; R = 2689.989305692122; S = 0.04700330116711059; T = 0.04700330116711057 }
    36 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 19.609342; N = 58.82802600000001; O = 0.06946983056124746; P = 0.06946983056124745; Q = 91.74425440830535; R = 825.6982896747481; S = 0.01442776939693614; T = 0.01442776939693614 }
    37 = @{ E = 3.0; G = 4.678599333333334; H = 14.035798; I = 0.2076839583510432; J = 0.2076839583510431; K = 3.0; M = 16.01808066666667; N = 48.054242; O = 0.05674710297927013; P = 0.05674710297927013; Q = 74.94218152834623; R = 674.479633755116; S = 0.0117854629716891; T = 0.01178546297168909 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
